$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) column values are written as text, matching the source data format
# (values like "1.005" or "7.010" must stay as literal text, not be parsed as numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.637.73'
$ws.Range("E2").Value = '  +2.12%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.868.89'
$ws.Range("E3").Value = '  +2.02%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.13'
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4633'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3883'
$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07871'
$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9754'
$ws.Range("E10").Value = '  +1.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.95'
$ws.Range("E11").Value = '  +0.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.863.14'
$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.010'
$ws.Range("E13").Value = '  +1.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.702'
$ws.Range("E14").Value = '  +0.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06978'
$ws.Range("E15").Value = '  +3.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.04'
$ws.Range("E16").Value = '  +0.97%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001003'
$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.81'
$ws.Range("E19").Value = '  +1.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.654.20'
$ws.Range("E21").Value = '  +2.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.282'
$ws.Range("E22").Value = '  -0.59%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.118'
$ws.Range("E24").Value = '  +1.07%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.124.82'
$ws.Range("E25").Value = '  +1.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.75'
$ws.Range("E26").Value = '  -0.60%  '

$ws.Range("E27").Value = '  +0.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.782'
$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.986'
$ws.Range("E29").Value = '  +0.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.32'
$ws.Range("E30").Value = '  +1.58%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09351'
$ws.Range("E31").Value = '  +0.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9182'
$ws.Range("E32").Value = '  -2.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.270'
$ws.Range("E33").Value = '  -0.49%  '

$ws.Range("E34").Value = '  +1.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.334'
$ws.Range("E35").Value = '  +1.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05782'
$ws.Range("E36").Value = '  -1.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02108'
$ws.Range("E37").Value = '  -1.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.152'
$ws.Range("E38").Value = '  +0.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.763'
$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5629'
$ws.Range("E40").Value = '  +0.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1789'
$ws.Range("E41").Value = '  +1.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.755'
$ws.Range("E42").Value = '  -1.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07197'
$ws.Range("E43").Value = '  +2.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.78'
$ws.Range("E44").Value = '  +1.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5314'
$ws.Range("E45").Value = '  +0.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.144'
$ws.Range("E46").Value = '  +3.51%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.110'
$ws.Range("E47").Value = '  -1.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.829'
$ws.Range("E48").Value = '  -0.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.21'
$ws.Range("E49").Value = '  +0.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.413'
$ws.Range("E50").Value = '  +4.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.52%  '
